$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.272.12"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "'2.276.99"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'113.95"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "'303.22"
$ws.Range("E6").Value = "  +7.44%  "
$ws.Range("D7").Value = "'0.634"
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "'44.76"
$ws.Range("E10").Value = "  -4.28%  "
$ws.Range("D11").Value = "'0.0928"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "'55.19"
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").Value = "'8.92"
$ws.Range("E13").Value = "  -3.26%  "
$ws.Range("E14").Value = "  +19.42%  "
$ws.Range("D15").Value = "'0.105"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "'15.46"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "'2.618.11"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").Value = "'2.272.45"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").Value = "'43.199.48"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  +5.71%  "
$ws.Range("D22").Value = "'75.31"
$ws.Range("E22").Value = "  +4.21%  "
$ws.Range("D23").Value = "'3.54"
$ws.Range("E23").Value = "  +11.40%  "
$ws.Range("D24").Value = "'257.06"
$ws.Range("E24").Value = "  +10.86%  "
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("D26").Value = "'9.09"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("D27").Value = "'11.73"
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'38.40"
$ws.Range("E30").Value = "  -5.13%  "
$ws.Range("D31").Value = "'22.31"
$ws.Range("E31").Value = "  +5.36%  "
$ws.Range("D32").Value = "'175.35"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("D34").Value = "'0.0898"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("D36").Value = "'5.11"
$ws.Range("E36").Value = "  +9.60%  "
$ws.Range("D37").Value = "'4.31"
$ws.Range("E37").Value = "  -6.48%  "
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").Value = "'2.47"
$ws.Range("E41").Value = "  -5.18%  "
$ws.Range("D42").Value = "'72.51"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'12.70"
$ws.Range("E45").Value = "  -5.82%  "
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").Value = "'5.65"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").Value = "'107.85"
$ws.Range("E48").Value = "  +6.83%  "
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").Value = "'8.77"
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("D51").Value = "'72.84"
$ws.Range("E51").Value = "  +3.89%  "
